# Appends a new data row (row 30) to both worksheets (G5368 and G4945),
# mirroring the existing row layout (country .. median_exceedance_pct).
# The new rows represent the 2025-10-29 / lead_time 3d forecast record
# that was generated when the analysis data was produced as xlsx instead of csv.

$wb = $excel.ActiveWorkbook

function Set-RowData {
    param($ws, $row, $values)

    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 1
        $cell = $ws.Cells.Item($row, $col)
        $val = $values[$i]

        if ($val -is [string] -and $val -match '^\d{4}-\d{2}-\d{2}$') {
            # Keep ISO-looking date strings (e.g. forecast_date) stored as plain
            # text instead of letting Excel auto-convert them to date serials.
            $cell.NumberFormat = "@"
            $cell.Value2 = $val
            $cell.ClearFormats()
        }
        else {
            $cell.Value2 = $val
        }
    }
}

# --- Sheet "G5368" (Nia Pumping Station) ---
$ws1 = $wb.Worksheets.Item("G5368")
$row1Values = @(
    "Philippines",
    "philippines",
    "Agusan River Basin",
    "agusan",
    "Nia Pumping Station",
    "G5368",
    "primary",
    "2025-10-29",
    3,
    8.874999999999865,
    125.5749999999995,
    5,
    4709.973879596918,
    "LOW",
    3193.342710267902,
    4709.973879596918,
    50,
    0,
    0,
    1034.98828125,
    1097.994262695312,
    764.4375,
    1826.0546875,
    920.83203125,
    1198.337890625,
    $false,
    -78.02560464860635
)
Set-RowData $ws1 30 $row1Values

# --- Sheet "G4945" (Talacogon Municipal Hall) ---
$ws2 = $wb.Worksheets.Item("G4945")
$row2Values = @(
    "Philippines",
    "philippines",
    "Agusan River Basin",
    "agusan",
    "Talacogon Municipal Hall",
    "G4945",
    "secondary",
    "2025-10-29",
    3,
    8.424999999999859,
    125.7749999999995,
    5,
    3363.250778297076,
    "LOW",
    2342.691130371584,
    3363.250778297076,
    50,
    0,
    0,
    684.63671875,
    736.4121704101562,
    450.953125,
    1507.0390625,
    581.0234375,
    782.27734375,
    $false,
    -79.64360186377021
)
Set-RowData $ws2 30 $row2Values
